# Reorders the account-statement worker rows so each worker's periods are
# grouped together (Gustavo Adolfo Galeano Oviedo first, most recent period
# 2112 down to 2106; then Luisa Edilma Rey Pabon, 2112 down to 2106) and
# refreshes the overdue/base-salary values to match the new source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2112", 21333, 1434505),
    @("CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2111", 40000, 1434505),
    @("CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2110", 57380, 1434505),
    @("CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2109", 57380, 1434505),
    @("CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2108", 57380, 1434505),
    @("CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2107", 57380, 1434505),
    @("CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2106", 40000, 1434505),
    @("CC", "37863609", "LUISA EDILMA REY PABON",        "2112", 32000, 1500000),
    @("CC", "37863609", "LUISA EDILMA REY PABON",        "2111", 60000, 1500000),
    @("CC", "37863609", "LUISA EDILMA REY PABON",        "2110", 60000, 1500000),
    @("CC", "37863609", "LUISA EDILMA REY PABON",        "2109", 60000, 1500000),
    @("CC", "37863609", "LUISA EDILMA REY PABON",        "2108", 60000, 1500000),
    @("CC", "37863609", "LUISA EDILMA REY PABON",        "2107", 60000, 1500000),
    @("CC", "37863609", "LUISA EDILMA REY PABON",        "2106", 60000, 1500000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
}
